$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(-2, -2, 5, 1, 5, -2, -3, 6, -4, 3, 3, 1, 7, -1, 4, 6, 4, -3, 0, 3, 1, -4, 0, 1, -1, -5, 0, -6, 1, -2, -1, 0, 0, 4, 2)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $values[$i]
}
